$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.26020348072052
$ws.Range("B1").Value = 1.628110766410828
$ws.Range("C1").Value = 2.340687036514282
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 1.358235955238342
